$d = $word.ActiveDocument

$d.Content.Find.Execute("392×9=3528", $true, $false, $false, $false, $false, $true, 1, $false, "471×9=4239", 2)
$d.Content.Find.Execute("174×5=870", $true, $false, $false, $false, $false, $true, 1, $false, "182×4=728", 2)
$d.Content.Find.Execute("424×3=1272", $true, $false, $false, $false, $false, $true, 1, $false, "676×8=5408", 2)
$d.Content.Find.Execute("802×6=4812", $true, $false, $false, $false, $false, $true, 1, $false, "396×5=1980", 2)
$d.Content.Find.Execute("514×6=3084", $true, $false, $false, $false, $false, $true, 1, $false, "401×5=2005", 2)
$d.Content.Find.Execute("904×2=1808", $true, $false, $false, $false, $false, $true, 1, $false, "169×2=338", 2)
$d.Content.Find.Execute("518×9=4662", $true, $false, $false, $false, $false, $true, 1, $false, "303×2=606", 2)
$d.Content.Find.Execute("139×5=695", $true, $false, $false, $false, $false, $true, 1, $false, "971×4=3884", 2)
$d.Content.Find.Execute("431×2=862", $true, $false, $false, $false, $false, $true, 1, $false, "113×4=452", 2)
$d.Content.Find.Execute("299×2=598", $true, $false, $false, $false, $false, $true, 1, $false, "665×8=5320", 2)
$d.Content.Find.Execute("333×4=1332", $true, $false, $false, $false, $false, $true, 1, $false, "885×8=7080", 2)
$d.Content.Find.Execute("422×2=844", $true, $false, $false, $false, $false, $true, 1, $false, "600×2=1200", 2)
$d.Content.Find.Execute("115×9=1035", $true, $false, $false, $false, $false, $true, 1, $false, "441×3=1323", 2)
$d.Content.Find.Execute("338×6=2028", $true, $false, $false, $false, $false, $true, 1, $false, "877×7=6139", 2)
$d.Content.Find.Execute("466×3=1398", $true, $false, $false, $false, $false, $true, 1, $false, "179×4=716", 2)
$d.Content.Find.Execute("985×5=4925", $true, $false, $false, $false, $false, $true, 1, $false, "145×4=580", 2)
$d.Content.Find.Execute("975×5=4875", $true, $false, $false, $false, $false, $true, 1, $false, "845×8=6760", 2)
$d.Content.Find.Execute("298×3=894", $true, $false, $false, $false, $false, $true, 1, $false, "944×5=4720", 2)
$d.Content.Find.Execute("238×9=2142", $true, $false, $false, $false, $false, $true, 1, $false, "743×8=5944", 2)
$d.Content.Find.Execute("506×5=2530", $true, $false, $false, $false, $false, $true, 1, $false, "180×7=1260", 2)
$d.Content.Find.Execute("701×8=5608", $true, $false, $false, $false, $false, $true, 1, $false, "565×5=2825", 2)
$d.Content.Find.Execute("204×9=1836", $true, $false, $false, $false, $false, $true, 1, $false, "215×9=1935", 2)
$d.Content.Find.Execute("437×9=3933", $true, $false, $false, $false, $false, $true, 1, $false, "461×6=2766", 2)
$d.Content.Find.Execute("861×4=3444", $true, $false, $false, $false, $false, $true, 1, $false, "558×6=3348", 2)
$d.Content.Find.Execute("985×2=1970", $true, $false, $false, $false, $false, $true, 1, $false, "343×6=2058", 2)
